# "started with data shadows" - slide 3 (Deprecate Compound Sets) rework:
#  - nudge the 5 connector/label shapes to their new positions/sizes
#  - add a "Add dcsu library." line ahead of the runtime-library bullets
#  - fix a typo: "main model" -> "To Main model"
#  - add a "Remove dcsu library" line after the compound-set cleanup bullets

# Convert OOXML EMU -> points for Shape.Left/Top/Width/Height. A tiny
# epsilon is added because the COM property setters round-trip the value
# through a single-precision float before re-deriving EMUs on save, which
# otherwise occasionally rounds the stored EMU one unit low.
function EMU([double]$emu) {
    return ($emu / 12700.0) + 0.00004
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# --- Rectangle: Rounded Corners 3 (id 4) -> reposition only ---
$shRect3 = $s.Shapes.Item(3)
$shRect3.Left = EMU 2667895
$shRect3.Top  = EMU 151473

# --- Oval 4 (id 5) -> reposition/resize + text edits ---
$shOval4 = $s.Shapes.Item(4)
$shOval4.Left   = EMU 1301672
$shOval4.Top    = EMU 822533
$shOval4.Width  = EMU 9143999
$shOval4.Height = EMU 2060415

$tr4 = $shOval4.TextFrame.TextRange
$tr4.Text = ""
$tr4.InsertAfter("Add ") | Out-Null
$tr4.InsertAfter("dcsu") | Out-Null
$tr4.InsertAfter(" library.`r") | Out-Null
$tr4.InsertAfter("Create Mapping Set Runtime Library with`r") | Out-Null
$tr4.InsertAfter("1. Mapping set, 2 data shadow, 3 Mapping procedures `r") | Out-Null
$tr4.InsertAfter("Create Content Type with mapping set and without compound`r") | Out-Null
$tr4.InsertAfter("Manually move the Mapping Set To Main model`r") | Out-Null
$tr4.InsertAfter("For each case: Copy data to mapping data in Runtime Library") | Out-Null

# --- Rectangle: Rounded Corners 5 (id 6) -> reposition only ---
$shRect5 = $s.Shapes.Item(5)
$shRect5.Left = EMU 2463501
$shRect5.Top  = EMU 3031473

# --- Oval 6 (id 7) -> reposition/resize + text edits ---
$shOval6 = $s.Shapes.Item(6)
$shOval6.Left   = EMU 345349
$shOval6.Top    = EMU 3742884
$shOval6.Width  = EMU 11501302
$shOval6.Height = EMU 1904879

$tr7 = $shOval6.TextFrame.TextRange
$tr7.Text = ""
$tr7.InsertAfter("Manually, move the indices from the compound sets to mapping sets`r") | Out-Null
$tr7.InsertAfter("Manually, remove the compound sets`r") | Out-Null
$tr7.InsertAfter("Manually, let the model builder fix compilation errors in the model`r") | Out-Null
$tr7.InsertAfter("For each case, copy the data back to the original identifiers`r") | Out-Null
$tr7.InsertAfter("Remove ") | Out-Null
$tr7.InsertAfter("dcsu") | Out-Null
$tr7.InsertAfter(" library") | Out-Null

# --- Rectangle: Rounded Corners 7 (id 8) -> reposition only ---
$shRect7 = $s.Shapes.Item(7)
$shRect7.Left = EMU 1992258
$shRect7.Top  = EMU 5741996
